$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; Excel shifts rows 12:66 down to 13:67
# and inherits the number formatting of the row above (so the date cell
# in column D keeps its date style automatically).
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new market-report entry.
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = "2022-09-27"
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112031
$ws.Range("G12").Value = "Poroto verde"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 1300
$ws.Range("L12").Value = 1400
$ws.Range("M12").Value = 1350
$ws.Range("N12").Value = "$/kilo"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 1350
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"
